$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New cumulative COVID death data rows (58-64), continuing the existing table.
$rows = @(
    @{ Row = 58; Date = 44175; B = 1122; C = 290; D = 1412 },
    @{ Row = 59; Date = 44176; B = 1148; C = 293; D = 1441 },
    @{ Row = 60; Date = 44177; B = 1175; C = 303; D = 1478 },
    @{ Row = 61; Date = 44178; B = 1205; C = 313; D = 1518 },
    @{ Row = 62; Date = 44179; B = 1251; C = 334; D = 1585 },
    @{ Row = 63; Date = 44180; B = 1309; C = 350; D = 1659 },
    @{ Row = 64; Date = 44181; B = 1378; C = 365; D = 1743 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Match the existing date-column formatting by copying it down from the
    # last populated row instead of minting a brand-new number format.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value2 = $r.Date

    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
}

# Update selection like the source workbook does after the edit.
$ws.Range("D51").Select()

$wb.Save()
